# Update for north-macedonia 1-mfl 2023-2024: reorder match rows 3-45 and add new row 46
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 is brand new - copy formatting (styles) from row 45 before writing values
$ws.Range("A45:V45").Copy()
$ws.Range("A46:V46").PasteSpecial(-4122)

# Target data for columns F:V across rows 3-46 (A:E already hold the correct
# index/metadata/date values and do not need to change, except row 46 below)
$rowsData = @(
  @{ Row=3; F="Tikves"; G=1; H="Makedonija GP"; I=0; J=2.24; K="05/08/2023 05:13"; L=2.34; M="06/08/2023 16:51"; N=2.84; O="05/08/2023 05:13"; P=3.1; Q="06/08/2023 16:12"; R=2.87; S="05/08/2023 05:13"; T=2.82; U="06/08/2023 16:51"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/tikves-makedonija-gp/xreALDzo/" },
  @{ Row=4; F="Struga"; G=2; H="KF Gostivar"; I=0; J=1.2; K="06/08/2023 11:43"; L=1.42; M="06/08/2023 16:59"; N=5.7; O="06/08/2023 11:43"; P=4.1; Q="06/08/2023 16:59"; R=10.45; S="06/08/2023 11:43"; T=6.35; U="06/08/2023 16:59"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/struga-kf-gostivar/0M3rQgST/" },
  @{ Row=5; F="Brera Strumica"; G=1; H="Vardar"; I=0; J=1.85; K="06/08/2023 11:43"; L=1.64; M="06/08/2023 16:19"; N=3.02; O="06/08/2023 11:43"; P=3.19; Q="06/08/2023 16:19"; R=4.15; S="06/08/2023 11:43"; T=5.47; U="06/08/2023 16:19"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/brera-strumica-vardar/S85zSFcH/" },
  @{ Row=6; F="Bregalnica Stip"; G=2; H="Sileks"; I=1; J=2.11; K="05/08/2023 05:12"; L=2.15; M="06/08/2023 16:08"; N=2.81; O="05/08/2023 05:12"; P=2.95; Q="06/08/2023 16:08"; R=3.12; S="05/08/2023 05:12"; T=3.32; U="06/08/2023 16:08"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/bregalnica-stip-sileks/vmZhpcDo/" },
  @{ Row=7; F="Rabotnicki"; G=3; H="Shkupi"; I=1; J=4.21; K="06/08/2023 05:12"; L=5; M="07/08/2023 16:55"; N=3.12; O="06/08/2023 05:12"; P=3.34; Q="07/08/2023 16:55"; R=1.69; S="06/08/2023 05:12"; T=1.65; U="07/08/2023 16:30"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/rabotnicki-shkupi/IwYdqHSi/" },
  @{ Row=8; F="Sileks"; G=2; H="Rabotnicki"; I=0; J=2.29; K="12/08/2023 05:12"; L=2.31; M="13/08/2023 16:57"; N=2.81; O="12/08/2023 05:12"; P=3.05; Q="13/08/2023 16:54"; R=2.82; S="12/08/2023 05:12"; T=2.56; U="13/08/2023 16:57"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/sileks-rabotnicki/j5o4sec4/" },
  @{ Row=9; F="Brera Strumica"; G=2; H="Voska Sport"; I=0; J=1.66; K="13/08/2023 10:38"; L=1.79; M="13/08/2023 15:40"; N=3.39; O="13/08/2023 10:38"; P=3.36; Q="13/08/2023 15:40"; R=4.63; S="13/08/2023 10:38"; T=4.03; U="13/08/2023 15:40"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/brera-strumica-voska-sport/KpwLwDdT/" },
  @{ Row=10; F="Vardar"; G=0; H="Shkupi"; I=1; J=3.95; K="13/08/2023 10:37"; L=5.57; M="13/08/2023 16:31"; N=3.15; O="13/08/2023 10:37"; P=3.41; Q="13/08/2023 16:31"; R=1.85; S="13/08/2023 10:37"; T=1.58; U="13/08/2023 16:31"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/vardar-shkupi/CEn0ryrb/" },
  @{ Row=11; F="Makedonija GP"; G=1; H="Bregalnica Stip"; I=0; J=2.1; K="12/08/2023 05:12"; L=1.98; M="13/08/2023 16:45"; N=2.88; O="12/08/2023 05:12"; P=2.99; Q="13/08/2023 16:12"; R=3.08; S="12/08/2023 05:12"; T=3.78; U="13/08/2023 16:45"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/makedonija-gp-bregalnica-stip/YTz9tFCA/" },
  @{ Row=12; F="KF Gostivar"; G=0; H="Tikves"; I=3; J=3.13; K="13/08/2023 10:38"; L=2.05; M="13/08/2023 16:59"; N=2.94; O="13/08/2023 10:38"; P=3.04; Q="13/08/2023 16:59"; R=2.21; S="13/08/2023 10:38"; T=3.45; U="13/08/2023 16:59"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/kf-gostivar-tikves/nXvDuZRG/" },
  @{ Row=13; F="Voska Sport"; G=3; H="Vardar"; I=2; J=1.74; K="20/08/2023 08:55"; L=1.95; M="20/08/2023 16:22"; N=3.33; O="20/08/2023 08:55"; P=3.15; Q="20/08/2023 16:22"; R=4.24; S="20/08/2023 08:55"; T=3.62; U="20/08/2023 16:22"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/voska-sport-vardar/0vLmZW4p/" },
  @{ Row=14; F="Shkupi"; G=1; H="Sileks"; I=1; J=1.31; K="19/08/2023 05:12"; L=1.4; M="20/08/2023 16:58"; N=4.13; O="19/08/2023 05:12"; P=4.3; Q="20/08/2023 16:58"; R=7.01; S="19/08/2023 05:12"; T=5.69; U="20/08/2023 16:58"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/shkupi-sileks/QNanfYtA/" },
  @{ Row=15; F="Rabotnicki"; G=1; H="Makedonija GP"; I=0; J=1.96; K="19/08/2023 05:12"; L=1.89; M="20/08/2023 16:59"; N=2.97; O="19/08/2023 05:12"; P=3.12; Q="20/08/2023 16:59"; R=3.33; S="19/08/2023 05:12"; T=3.88; U="20/08/2023 16:59"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/rabotnicki-makedonija-gp/nX0reER3/" },
  @{ Row=16; F="Bregalnica Stip"; G=2; H="KF Gostivar"; I=1; J=1.67; K="20/08/2023 08:55"; L=1.29; M="20/08/2023 16:59"; N=3.49; O="20/08/2023 08:55"; P=4.55; Q="20/08/2023 16:59"; R=4.41; S="20/08/2023 08:55"; T=6.27; U="20/08/2023 16:59"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/bregalnica-stip-kf-gostivar/48BwdfCc/" },
  @{ Row=17; F="Tikves"; G=1; H="Shkendija"; I=2; J=4.42; K="19/08/2023 05:12"; L=3.91; M="20/08/2023 16:41"; N=3.3; O="19/08/2023 05:12"; P=3.28; Q="20/08/2023 16:41"; R=1.61; S="19/08/2023 05:12"; T=1.83; U="20/08/2023 16:41"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/tikves-shkendija-tetovo/hY7Zdzci/" },
  @{ Row=18; F="Struga"; G=1; H="Brera Strumica"; I=2; J=1.48; K="20/08/2023 05:12"; L=2.23; M="21/08/2023 16:51"; N=3.47; O="20/08/2023 05:12"; P=2.97; Q="21/08/2023 16:47"; R=5.3; S="20/08/2023 05:12"; T=3.01; U="21/08/2023 16:51"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/struga-brera-strumica/vVJiYjKj/" },
  @{ Row=19; F="KF Gostivar"; G=2; H="Rabotnicki"; I=0; J=2.63; K="23/08/2023 11:12"; L=3.52; M="23/08/2023 16:26"; N=3.03; O="23/08/2023 11:12"; P=3.25; Q="23/08/2023 16:26"; R=2.5; S="23/08/2023 11:12"; T=1.79; U="23/08/2023 16:26"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/kf-gostivar-rabotnicki/GChaiWQS/" },
  @{ Row=20; F="Vardar"; G=2; H="Sileks"; I=2; J=3.32; K="23/08/2023 11:12"; L=3.1; M="23/08/2023 15:31"; N=2.89; O="23/08/2023 11:12"; P=2.82; Q="23/08/2023 15:31"; R=2.15; S="23/08/2023 11:12"; T=2.34; U="23/08/2023 15:31"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/vardar-sileks/fc4jghdG/" },
  @{ Row=21; F="Shkendija"; G=2; H="Bregalnica Stip"; I=0; J=1.33; K="22/08/2023 04:12"; L=1.38; M="23/08/2023 16:29"; N=4.05; O="22/08/2023 04:12"; P=4.14; Q="23/08/2023 16:29"; R=6.47; S="22/08/2023 04:12"; T=5.23; U="23/08/2023 16:29"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/shkendija-tetovo-bregalnica-stip/hYqRnjZq/" },
  @{ Row=22; F="Makedonija GP"; G=0; H="Shkupi"; I=0; J=4.04; K="22/08/2023 04:12"; L=4.11; M="23/08/2023 15:56"; N=3.06; O="22/08/2023 04:12"; P=3.16; Q="23/08/2023 16:28"; R=1.74; S="22/08/2023 04:12"; T=1.83; U="23/08/2023 16:28"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/makedonija-gp-shkupi/Ym3fhCBM/" },
  @{ Row=23; F="Brera Strumica"; G=2; H="Tikves"; I=0; J=1.99; K="24/08/2023 16:12"; L=1.99; M="24/08/2023 16:12"; N=3.06; O="24/08/2023 16:12"; P=3.06; Q="24/08/2023 16:12"; R=3.23; S="24/08/2023 16:12"; T=3.23; U="24/08/2023 16:12"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/brera-strumica-tikves/neqVoAlj/" },
  @{ Row=24; F="Shkupi"; G=3; H="KF Gostivar"; I=1; J=1.27; K="27/08/2023 15:13"; L=1.3; M="27/08/2023 15:53"; N=4.94; O="27/08/2023 15:13"; P=4.64; Q="27/08/2023 15:53"; R=8.119999999999999; S="27/08/2023 15:13"; T=8.25; U="27/08/2023 15:53"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/shkupi-kf-gostivar/hh4jkqer/" },
  @{ Row=25; F="Rabotnicki"; G=1; H="Shkendija"; I=2; J=3.14; K="26/08/2023 04:12"; L=2.98; M="27/08/2023 15:59"; N=3.18; O="26/08/2023 04:12"; P=3.25; Q="27/08/2023 15:59"; R=1.98; S="26/08/2023 04:12"; T=1.96; U="27/08/2023 15:59"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/rabotnicki-shkendija-tetovo/CfjM6BQF/" },
  @{ Row=26; F="Sileks"; G=2; H="Makedonija GP"; I=2; J=2.05; K="26/08/2023 04:12"; L=2.12; M="27/08/2023 15:52"; N=2.84; O="26/08/2023 04:12"; P=2.8; Q="27/08/2023 15:52"; R=3.25; S="26/08/2023 04:12"; T=3.62; U="27/08/2023 15:52"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/sileks-makedonija-gp/Wr3fl3Al/" },
  @{ Row=27; F="Bregalnica Stip"; G=1; H="Brera Strumica"; I=0; J=2.48; K="26/08/2023 04:12"; L=2.32; M="27/08/2023 15:19"; N=2.78; O="26/08/2023 04:12"; P=3.01; Q="27/08/2023 15:33"; R=2.61; S="26/08/2023 04:12"; T=2.93; U="27/08/2023 15:19"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/bregalnica-stip-brera-strumica/binI7iB9/" },
  @{ Row=28; F="Voska Sport"; G=2; H="Bregalnica Stip"; I=0; J=2.54; K="02/09/2023 13:12"; L=2.02; M="02/09/2023 15:41"; N=2.81; O="02/09/2023 13:12"; P=3.13; Q="02/09/2023 15:41"; R=2.7; S="02/09/2023 13:12"; T=3.43; U="02/09/2023 15:41"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/voska-sport-bregalnica-stip/YL2h94nE/" },
  @{ Row=29; F="Vardar"; G=2; H="Makedonija GP"; I=1; J=2.31; K="02/09/2023 13:12"; L=2.69; M="02/09/2023 15:52"; N=2.93; O="02/09/2023 13:12"; P=2.97; Q="02/09/2023 15:49"; R=2.88; S="02/09/2023 13:12"; T=2.52; U="02/09/2023 15:52"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/vardar-makedonija-gp/0I7bmNPf/" },
  @{ Row=30; F="Shkendija"; G=0; H="Shkupi"; I=0; J=2.05; K="01/09/2023 04:12"; L=1.74; M="02/09/2023 15:58"; N=2.83; O="01/09/2023 04:12"; P=2.95; Q="02/09/2023 15:58"; R=3.26; S="01/09/2023 04:12"; T=4.2; U="02/09/2023 15:46"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/shkendija-tetovo-shkupi/MXe7o1f7/" },
  @{ Row=31; F="KF Gostivar"; G=1; H="Sileks"; I=2; J=2.67; K="01/09/2023 04:12"; L=2; M="02/09/2023 15:58"; N=2.78; O="01/09/2023 04:12"; P=3.11; Q="02/09/2023 15:58"; R=2.43; S="01/09/2023 04:12"; T=3.52; U="02/09/2023 15:58"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/kf-gostivar-sileks/v962nsu1/" },
  @{ Row=32; F="Brera Strumica"; G=0; H="Rabotnicki"; I=0; J=1.88; K="01/09/2023 04:12"; L=1.68; M="02/09/2023 15:49"; N=2.98; O="01/09/2023 04:12"; P=3.54; Q="02/09/2023 15:49"; R=3.57; S="01/09/2023 04:12"; T=4.38; U="02/09/2023 15:49"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/brera-strumica-rabotnicki/KIEApL9D/" },
  @{ Row=33; F="Struga"; G=3; H="Tikves"; I=1; J=1.45; K="03/09/2023 04:12"; L=1.47; M="04/09/2023 15:47"; N=3.68; O="03/09/2023 04:12"; P=3.66; Q="04/09/2023 15:47"; R=5.22; S="03/09/2023 04:12"; T=6.46; U="04/09/2023 15:47"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/struga-tikves/CC1d8O1K/" },
  @{ Row=34; F="Shkupi"; G=1; H="Brera Strumica"; I=0; J=1.56; K="15/09/2023 03:12"; L=1.61; M="16/09/2023 14:45"; N=3.24; O="15/09/2023 03:12"; P=3.45; Q="16/09/2023 14:45"; R=4.94; S="15/09/2023 03:12"; T=5.1; U="16/09/2023 14:45"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/shkupi-brera-strumica/vNoy1tff/" },
  @{ Row=35; F="Tikves"; G=1; H="Vardar"; I=0; J=1.73; K="16/09/2023 03:13"; L=1.99; M="17/09/2023 14:51"; N=3.13; O="16/09/2023 03:13"; P=3.05; Q="17/09/2023 14:51"; R=4; S="16/09/2023 03:13"; T=3.62; U="17/09/2023 14:51"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/tikves-vardar/Umd17rHQ/" },
  @{ Row=36; F="Sileks"; G=2; H="Shkendija"; I=1; J=4.23; K="16/09/2023 03:13"; L=3.86; M="17/09/2023 14:55"; N=3.14; O="16/09/2023 03:13"; P=2.61; Q="17/09/2023 14:55"; R=1.68; S="16/09/2023 03:13"; T=1.96; U="17/09/2023 14:55"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/sileks-shkendija-tetovo/Moyt0090/" },
  @{ Row=37; F="Rabotnicki"; G=2; H="Voska Sport"; I=0; J=1.63; K="17/09/2023 13:42"; L=2.24; M="17/09/2023 14:55"; N=3.41; O="17/09/2023 13:42"; P=3.09; Q="17/09/2023 14:55"; R=4.74; S="17/09/2023 13:42"; T=2.99; U="17/09/2023 14:55"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/rabotnicki-voska-sport/QJsX1Mvl/" },
  @{ Row=38; F="Makedonija GP"; G=3; H="KF Gostivar"; I=0; J=1.98; K="16/09/2023 03:13"; L=2.58; M="17/09/2023 14:56"; N=2.91; O="16/09/2023 03:13"; P=2.81; Q="17/09/2023 14:47"; R=3.33; S="16/09/2023 03:13"; T=2.67; U="17/09/2023 14:56"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/makedonija-gp-kf-gostivar/0xzpaKO6/" },
  @{ Row=39; F="Bregalnica Stip"; G=1; H="Struga"; I=3; J=3.67; K="16/09/2023 03:13"; L=3.58; M="17/09/2023 14:51"; N=2.98; O="16/09/2023 03:13"; P=3.13; Q="17/09/2023 14:51"; R=1.85; S="16/09/2023 03:13"; T=1.97; U="17/09/2023 14:51"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/bregalnica-stip-struga/WdrT22Ps/" },
  @{ Row=40; F="Brera Strumica"; G=1; H="Sileks"; I=1; J=1.9; K="23/09/2023 02:13"; L=1.9; M="24/09/2023 12:16"; N=2.93; O="23/09/2023 02:13"; P=3.06; Q="24/09/2023 13:04"; R=3.55; S="23/09/2023 02:13"; T=3.89; U="24/09/2023 12:16"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/brera-strumica-sileks/lbCzMsWJ/" },
  @{ Row=41; F="Vardar"; G=0; H="KF Gostivar"; I=2; J=2.32; K="23/09/2023 02:13"; L=2.31; M="24/09/2023 14:58"; N=2.81; O="23/09/2023 02:13"; P=2.79; Q="24/09/2023 14:58"; R=2.79; S="23/09/2023 02:13"; T=3.19; U="24/09/2023 14:58"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/vardar-kf-gostivar/82ZmbvvD/" },
  @{ Row=42; F="Tikves"; G=2; H="Bregalnica Stip"; I=1; J=2.17; K="23/09/2023 02:13"; L=2.02; M="24/09/2023 14:50"; N=2.89; O="23/09/2023 02:13"; P=2.81; Q="24/09/2023 14:50"; R=2.93; S="23/09/2023 02:13"; T=3.95; U="24/09/2023 14:50"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/tikves-bregalnica-stip/bX1AFu9m/" },
  @{ Row=43; F="Struga"; G=1; H="Rabotnicki"; I=0; J=1.44; K="23/09/2023 02:13"; L=1.63; M="24/09/2023 14:14"; N=3.65; O="23/09/2023 02:13"; P=3.43; Q="24/09/2023 14:14"; R=5.46; S="23/09/2023 02:13"; T=4.99; U="24/09/2023 14:14"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/struga-rabotnicki/Kt36GLgs/" },
  @{ Row=44; F="Shkendija"; G=1; H="Makedonija GP"; I=0; J=1.36; K="23/09/2023 02:13"; L=1.58; M="24/09/2023 14:45"; N=3.89; O="23/09/2023 02:13"; P=3.46; Q="24/09/2023 14:50"; R=6.28; S="23/09/2023 02:13"; T=5.48; U="24/09/2023 14:50"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/shkendija-tetovo-makedonija-gp/ngYicbgJ/" },
  @{ Row=45; F="Voska Sport"; G=2; H="Shkupi"; I=3; J=3.35; K="24/09/2023 12:13"; L=4.03; M="24/09/2023 14:32"; N=3.2; O="24/09/2023 12:13"; P=3.35; Q="24/09/2023 14:32"; R=2.02; S="24/09/2023 12:13"; T=1.79; U="24/09/2023 14:32"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/voska-sport-shkupi/21GvL1oQ/" },
  @{ Row=46; F="Voska Sport"; G=0; H="Struga"; I=1; J=3.69; K="26/09/2023 02:12"; L=3.08; M="27/09/2023 14:52"; N=3; O="26/09/2023 02:12"; P=3.14; Q="27/09/2023 14:52"; R=1.84; S="26/09/2023 02:12"; T=2.17; U="27/09/2023 14:52"; V="https://www.betexplorer.com/football/north-macedonia/1-mfl/voska-sport-struga/QorZpU3d/" }
)

foreach ($r in $rowsData) {
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
    $ws.Cells.Item($r.Row, 10).Value = $r.J
    $ws.Cells.Item($r.Row, 11).Value = $r.K
    $ws.Cells.Item($r.Row, 12).Value = $r.L
    $ws.Cells.Item($r.Row, 13).Value = $r.M
    $ws.Cells.Item($r.Row, 14).Value = $r.N
    $ws.Cells.Item($r.Row, 15).Value = $r.O
    $ws.Cells.Item($r.Row, 16).Value = $r.P
    $ws.Cells.Item($r.Row, 17).Value = $r.Q
    $ws.Cells.Item($r.Row, 18).Value = $r.R
    $ws.Cells.Item($r.Row, 19).Value = $r.S
    $ws.Cells.Item($r.Row, 20).Value = $r.T
    $ws.Cells.Item($r.Row, 21).Value = $r.U
    $ws.Cells.Item($r.Row, 22).Value = $r.V
}

# New row 46 needs its index (A) and shared metadata (B:E) populated too
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "north-macedonia"
$ws.Cells.Item(46, 3).Value = "1-mfl"
$ws.Cells.Item(46, 4).Value = "2023-2024"
$ws.Cells.Item(46, 5).Value = 45196.625
